$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.008135666666666668
$ws.Range("H2").Value = 0.024407
$ws.Range("I2").Value = 0.0001175588769867851
$ws.Range("J2").Value = 0.0001175588769867851
$ws.Range("M2").Value = 13.79081533333333
$ws.Range("N2").Value = 41.372446
$ws.Range("O2").Value = 0.2771682650956311
$ws.Range("P2").Value = 0.2771682650956311
$ws.Range("Q2").Value = 0.1121974766135556
$ws.Range("R2").Value = 1.009777289522
$ws.Range("S2").Value = 0.00003258358998101795
$ws.Range("T2").Value = 0.00003258358998101794

# Row 3
$ws.Range("G3").Value = 0.008135666666666668
$ws.Range("H3").Value = 0.024407
$ws.Range("I3").Value = 0.0001175588769867851
$ws.Range("J3").Value = 0.0001175588769867851
$ws.Range("O3").Value = 0.5287054547944754
$ws.Range("P3").Value = 0.5287054547944753
$ws.Range("Q3").Value = 0.2140195158320001
$ws.Range("R3").Value = 1.926175642488
$ws.Range("S3").Value = 0.00006215401952242602
$ws.Range("T3").Value = 0.000062154019522426

# Row 4
$ws.Range("G4").Value = 0.008135666666666668
$ws.Range("H4").Value = 0.024407
$ws.Range("I4").Value = 0.0001175588769867851
$ws.Range("J4").Value = 0.0001175588769867851
$ws.Range("M4").Value = 9.588979333333333
$ws.Range("N4").Value = 28.766938
$ws.Range("O4").Value = 0.1927196254621635
$ws.Range("P4").Value = 0.1927196254621635
$ws.Range("Q4").Value = 0.07801273952955556
$ws.Range("R4").Value = 0.7021146557660001
$ws.Range("S4").Value = 0.00002265590274264578
$ws.Range("T4").Value = 0.00002265590274264578

# Row 5
$ws.Range("G5").Value = 0.008135666666666668
$ws.Range("H5").Value = 0.024407
$ws.Range("I5").Value = 0.0001175588769867851
$ws.Range("J5").Value = 0.0001175588769867851
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.06998966666666667
$ws.Range("N5").Value = 0.209969
$ws.Range("O5").Value = 0.001406654647730148
$ws.Range("P5").Value = 0.001406654647730148
$ws.Range("Q5").Value = 0.0005694125981111112
$ws.Range("R5").Value = 0.005124713383000001
$ws.Range("S5").Value = 0.000000165364740695398
$ws.Range("T5").Value = 0.000000165364740695398

# Row 6
$ws.Range("I6").Value = 0.0001954433790098242
$ws.Range("J6").Value = 0.0001954433790098242
$ws.Range("M6").Value = 13.79081533333333
$ws.Range("N6").Value = 41.372446
$ws.Range("O6").Value = 0.2771682650956311
$ws.Range("P6").Value = 0.2771682650956311
$ws.Range("Q6").Value = 0.1865299712602222
$ws.Range("R6").Value = 1.678769741342
$ws.Range("S6").Value = 0.00005417070228458086
$ws.Range("T6").Value = 0.00005417070228458086

# Row 7
$ws.Range("I7").Value = 0.0001954433790098242
$ws.Range("J7").Value = 0.0001954433790098242
$ws.Range("O7").Value = 0.5287054547944754
$ws.Range("P7").Value = 0.5287054547944753
$ws.Range("Q7").Value = 0.3558106237520001
$ws.Range("S7").Value = 0.0001033319805859582
$ws.Range("T7").Value = 0.0001033319805859581

# Row 8
$ws.Range("I8").Value = 0.0001954433790098242
$ws.Range("J8").Value = 0.0001954433790098242
$ws.Range("M8").Value = 9.588979333333333
$ws.Range("N8").Value = 28.766938
$ws.Range("O8").Value = 0.1927196254621635
$ws.Range("P8").Value = 0.1927196254621635
$ws.Range("Q8").Value = 0.1296973381362222
$ws.Range("R8").Value = 1.167276043226
$ws.Range("S8").Value = 0.00003766577480183299
$ws.Range("T8").Value = 0.00003766577480183298

# Row 9
$ws.Range("I9").Value = 0.0001954433790098242
$ws.Range("J9").Value = 0.0001954433790098242
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.06998966666666667
$ws.Range("N9").Value = 0.209969
$ws.Range("O9").Value = 0.001406654647730148
$ws.Range("P9").Value = 0.001406654647730148
$ws.Range("Q9").Value = 0.0009466569014444445
$ws.Range("R9").Value = 0.008519912113000001
$ws.Range("S9").Value = 0.0000002749213374522541
$ws.Range("T9").Value = 0.0000002749213374522541

# Row 10
$ws.Range("G10").Value = 69.18337766666667
$ws.Range("H10").Value = 207.550133
$ws.Range("I10").Value = 0.9996869977440035
$ws.Range("J10").Value = 0.9996869977440034
$ws.Range("M10").Value = 13.79081533333333
$ws.Range("N10").Value = 41.372446
$ws.Range("O10").Value = 0.2771682650956311
$ws.Range("P10").Value = 0.2771682650956311
$ws.Range("Q10").Value = 954.0951855372576
$ws.Range("R10").Value = 8586.856669835319
$ws.Range("S10").Value = 0.2770815108033655
$ws.Range("T10").Value = 0.2770815108033655

# Row 11
$ws.Range("G11").Value = 69.18337766666667
$ws.Range("H11").Value = 207.550133
$ws.Range("I11").Value = 0.9996869977440035
$ws.Range("J11").Value = 0.9996869977440034
$ws.Range("O11").Value = 0.5287054547944754
$ws.Range("P11").Value = 0.5287054547944753
$ws.Range("Q11").Value = 1819.960625047208
$ws.Range("R11").Value = 16379.64562542488
$ws.Range("S11").Value = 0.528539968794367
$ws.Range("T11").Value = 0.5285399687943669

# Row 12
$ws.Range("G12").Value = 69.18337766666667
$ws.Range("H12").Value = 207.550133
$ws.Range("I12").Value = 0.9996869977440035
$ws.Range("J12").Value = 0.9996869977440034
$ws.Range("M12").Value = 9.588979333333333
$ws.Range("N12").Value = 28.766938
$ws.Range("O12").Value = 0.1927196254621635
$ws.Range("P12").Value = 0.1927196254621635
$ws.Range("Q12").Value = 663.3979786558616
$ws.Range("R12").Value = 5970.581807902754
$ws.Range("S12").Value = 0.192659303784619
$ws.Range("T12").Value = 0.192659303784619

# Row 13
$ws.Range("G13").Value = 69.18337766666667
$ws.Range("H13").Value = 207.550133
$ws.Range("I13").Value = 0.9996869977440035
$ws.Range("J13").Value = 0.9996869977440034
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.06998966666666667
$ws.Range("N13").Value = 0.209969
$ws.Range("O13").Value = 0.001406654647730148
$ws.Range("P13").Value = 0.001406654647730148
$ws.Range("Q13").Value = 4.842121541764111
$ws.Range("R13").Value = 43.57909387587701
$ws.Range("S13").Value = 0.001406214361652
$ws.Range("T13").Value = 0.001406214361652
